$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.373.40"
$ws.Range("E2").Value = "  -0.56%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.721.54"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.28"
$ws.Range("E5").Value = "  -0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4878"
$ws.Range("E7").Value = "  +1.85%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2611"
$ws.Range("E8").Value = "  -2.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06209"
$ws.Range("E9").Value = "  -0.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.722.46"
$ws.Range("E10").Value = "  -0.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07014"
$ws.Range("E11").Value = "  -2.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.47"
$ws.Range("E12").Value = "  -1.29%  "

$ws.Range("E13").Value = "  +0.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5964"
$ws.Range("E14").Value = "  -2.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.20"
$ws.Range("E15").Value = "  +0.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.392.74"
$ws.Range("E17").Value = "  -0.49%  "

$ws.Range("E18").Value = "  -0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007230"
$ws.Range("E19").Value = "  +3.80%  "

$ws.Range("E20").Value = "  -2.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.948.08"
$ws.Range("E21").Value = "  -0.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.477"
$ws.Range("E22").Value = "  -1.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.556"
$ws.Range("E23").Value = "  -3.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.165"
$ws.Range("E24").Value = "  -2.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.33"
$ws.Range("E25").Value = "  +0.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.23"
$ws.Range("E26").Value = "  -0.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.418"
$ws.Range("E27").Value = "  +1.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "107.17"
$ws.Range("E28").Value = "  +0.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.714"
$ws.Range("E29").Value = "  -4.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.958"
$ws.Range("E30").Value = "  -0.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07951"
$ws.Range("E31").Value = "  -0.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.676"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04526"
$ws.Range("E33").Value = "  -1.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.614"
$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9947"
$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6209"
$ws.Range("E36").Value = "  -1.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9064"
$ws.Range("E37").Value = "  -0.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.980"
$ws.Range("E38").Value = "  -5.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.386"
$ws.Range("E39").Value = "  -0.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.0000"
$ws.Range("E40").Value = "  -0.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01486"
$ws.Range("E41").Value = "  -0.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.31"
$ws.Range("E42").Value = "  -4.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.379"
$ws.Range("E43").Value = "  -3.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3847"
$ws.Range("E44").Value = "  -0.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.728"
$ws.Range("E45").Value = "  -3.82%  "

$ws.Range("E46").Value = "  -2.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05353"
$ws.Range("E47").Value = "  +0.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.06"
$ws.Range("E48").Value = "  -2.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.669"
$ws.Range("E49").Value = "  -2.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.240"
$ws.Range("E50").Value = "  -1.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.00"
$ws.Range("E51").Value = "  -0.38%  "
